$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (2021-08-11)
$ws.Range("B2").Value = 0.01253208636536152
$ws.Range("C2").Value = 0.04103571897497393
$ws.Range("D2").Value = 0.7210945179870265
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 1.308048281929061

# Row 3 (2021-08-10)
$ws.Range("B3").Value = 0.04172184405617529
$ws.Range("C3").Value = 0.3048912486333797
$ws.Range("D3").Value = 3993.344853322108
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 4007.555312885598

# Row 4 (2021-08-07)
$ws.Range("B4").Value = 0.00000009142958989905026
$ws.Range("C4").Value = 0.04103571897497393
$ws.Range("D4").Value = 0.1496068669990043
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 0.7240286360052668

# Row 5 (2021-07-22)
$ws.Range("B5").Value = 1.445647641019636
$ws.Range("C5").Value = 0.3048912486333797
$ws.Range("D5").Value = 189.6080260415259
$ws.Range("E5").Value = 13.86384647080068
$ws.Range("G5").Value = 205.2224114019796
